$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The old sheet had a 2-row header (units split across row1/row2).
# The new layout uses a single header row with new columns (idx, idx2, Name,
# Date Start, Date End) prepended and renamed unit headers. Remove the old
# second header row - this shifts all the data rows up by one and leaves
# row 1 as a single blank header row ready to be rewritten.
$ws.Rows.Item(2).Delete()

# Rewrite the (now single) header row with the new column headers.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# The new leading columns (idx/idx2/Name/Date Start/Date End) use the
# sheet's plain default formatting - clear any inherited formatting from
# the cells that used to live there (e.g. old E1 carried the old unit
# header's style).
$ws.Range("A1:E1").ClearFormats()

# Match the font formatting used by the rest of the header cells (Arial 9,
# same as the other unit-header cells in the sheet).
$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").Font.Name = "Arial"

# Match the selection state recorded after the edit.
$ws.Range("A2:K2").Select()
